$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 - header/count values
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2
$ws.Range("B2").Value = 15.350000000000001
$ws.Range("C2").Value = 7.3000000000000007
$ws.Range("D2").Value = 13.25
$ws.Range("E2").Value = 16.350000000000001

# Row 3
$ws.Range("B3").Value = 6.8500000000000005
$ws.Range("C3").Value = 17.05
$ws.Range("D3").Value = 23.8
$ws.Range("E3").Value = 13.850000000000001

# Update selection to match the edited range
$ws.Range("B1:E3").Select()
